# Auto-generated from the cryptos.xlsx OOXML diff.
# Updates Price (D) / Volume(1h) (E) text cells, and for the two rows
# whose ranking swapped (39<->40, 49<->50) also Coin (B) and Link (C).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.341.92"
$ws.Range("E2").Value = "  +3.55%  "
$ws.Range("D3").Value = "2.616.08"
$ws.Range("E3").Value = "  +1.65%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'570.44"
$ws.Range("E5").Value = "  +6.09%  "
$ws.Range("D6").Value = "'145.68"
$ws.Range("E6").Value = "  +2.19%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  +4.19%  "
$ws.Range("D9").Value = "2.632.58"
$ws.Range("E9").Value = "  +2.38%  "
$ws.Range("D10").Value = "'6.71"
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("E11").Value = "  +5.18%  "
$ws.Range("E12").Value = "  +11.11%  "
$ws.Range("D13").Value = "'0.344"
$ws.Range("E13").Value = "  +3.91%  "
$ws.Range("D14").Value = "3.084.66"
$ws.Range("E14").Value = "  +1.86%  "
$ws.Range("D15").Value = "60.347.74"
$ws.Range("E15").Value = "  +3.72%  "
$ws.Range("D16").Value = "'22.12"
$ws.Range("E16").Value = "  +7.75%  "
$ws.Range("D17").Value = "'0.0000138"
$ws.Range("E17").Value = "  +4.92%  "
$ws.Range("D18").Value = "2.636.72"
$ws.Range("E18").Value = "  +3.47%  "
$ws.Range("E19").Value = "  +2.31%  "
$ws.Range("D20").Value = "'342.21"
$ws.Range("E20").Value = "  +2.66%  "
$ws.Range("E21").Value = "  +4.13%  "
$ws.Range("D22").Value = "'6.36"
$ws.Range("E22").Value = "  +4.17%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "'65.73"
$ws.Range("E24").Value = "  -0.88%  "
$ws.Range("D25").Value = "'0.451"
$ws.Range("E25").Value = "  +7.91%  "
$ws.Range("E26").Value = "  +4.38%  "
$ws.Range("D27").Value = "'0.997"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").Value = "'7.37"
$ws.Range("E28").Value = "  +5.44%  "
$ws.Range("D29").Value = "0.0₃0797"
$ws.Range("E29").Value = "  +9.50%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  +4.81%  "
$ws.Range("E32").Value = "  +4.14%  "
$ws.Range("D33").Value = "'160.46"
$ws.Range("E33").Value = "  +3.07%  "
$ws.Range("D34").Value = "'19.17"
$ws.Range("E34").Value = "  +1.76%  "
$ws.Range("D35").Value = "'4.13"
$ws.Range("E35").Value = "  +6.75%  "
$ws.Range("D36").Value = "'0.902"
$ws.Range("E36").Value = "  +10.88%  "
$ws.Range("E37").Value = "  +6.16%  "
$ws.Range("D38").Value = "'0.883"
$ws.Range("E38").Value = "  +5.00%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'1.52"
$ws.Range("E39").Value = "  +7.81%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'37.52"
$ws.Range("E40").Value = "  +1.86%  "
$ws.Range("D41").Value = "'298.65"
$ws.Range("E41").Value = "  +7.22%  "
$ws.Range("D42").Value = "'3.66"
$ws.Range("E42").Value = "  +2.74%  "
$ws.Range("D43").Value = "'0.996"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").Value = "'0.0984"
$ws.Range("E44").Value = "  +4.93%  "
$ws.Range("D45").Value = "'0.601"
$ws.Range("E45").Value = "  +2.66%  "
$ws.Range("D46").Value = "'0.0544"
$ws.Range("E46").Value = "  +3.01%  "
$ws.Range("D47").Value = "'19.36"
$ws.Range("E47").Value = "  +5.47%  "
$ws.Range("D48").Value = "'10.69"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'126.42"
$ws.Range("E49").Value = "  +17.02%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0236"
$ws.Range("E50").Value = "  +4.58%  "
$ws.Range("D51").Value = "'18.74"
$ws.Range("E51").Value = "  +5.82%  "
